$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The post about "「アディダス」" (adidas) that used to live in row 32 was
# removed. Deleting the entire row shifts every subsequent row up by one,
# which matches the diff (row 33 "「コカコーラ・ゼロ」" becomes the new
# row 32, ..., the old row 206 becomes row 205) and shrinks the sheet's
# used range from A1:C206 to A1:C205.
$ws.Rows.Item(32).Delete()
